# Split single-run paragraphs into multiple runs joined by manual line breaks
# (<w:br/>) at the sentence/topic boundaries, per the source diff.
$d = $word.ActiveDocument

function Split-WithLineBreaks {
    param(
        [string]$FindText,
        [string]$ReplaceText
    )
    $rng = $d.Content
    $found = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        # Assigning Range.Text (rather than Find.Execute's Replacement arg) avoids
        # Word's smart-quote autocorrect mangling straight quotes in the source text.
        $rng.Text = $ReplaceText
    }
    return $found
}

# Edit 1
$ok1 = Split-WithLineBreaks "Ensinar a linguagem gráfica normalizada internacionalmente para representação de máquinas e equipamentos que integram os processos de engenharia.Desenvolver o raciocínio espacial e a criatividade de representação." "Ensinar a linguagem gráfica normalizada internacionalmente para representação de máquinas e equipamentos que integram os processos de engenharia.`vDesenvolver o raciocínio espacial e a criatividade de representação."
Write-Output "Edit 1: $ok1"

# Edit 2
$ok2 = Split-WithLineBreaks "-Introdução-Teoria Elementar do Desenho Projetivo-Projeções Ortogonais pelo 1º Diedro-Projeções Ortogonais pelo 3º Diedro-Leitura e Interpretação de Desenhos-Escalas-Desenhos com Instrumentos-Cortes e Representações Convencionais-Projeções Auxiliares-Cotação-Desenhos de Conjuntos e Detalhes-Aplicação de Tolerâncias e Ajustes -Símbolos de Acabamento Superficial-Desenho de Elementos de Máquina-Desenho de Equipamentos e Acessórios" "-Introdução`v-Teoria Elementar do Desenho Projetivo`v-Projeções Ortogonais pelo 1º Diedro`v-Projeções Ortogonais pelo 3º Diedro`v-Leitura e Interpretação de Desenhos`v-Escalas`v-Desenhos com Instrumentos`v-Cortes e Representações Convencionais`v-Projeções Auxiliares`v-Cotação`v-Desenhos de Conjuntos e Detalhes`v-Aplicação de Tolerâncias e Ajustes `v-Símbolos de Acabamento Superficial`v-Desenho de Elementos de Máquina`v-Desenho de Equipamentos e Acessórios"
Write-Output "Edit 2: $ok2"

# Edit 3
$ok3 = Split-WithLineBreaks "1 - INTRODUÇÃOApresentação e definição da disciplina, destacando a importância do desenho na engenharia; Normas ABNT e ISO.2 - TEORIA ELEMENTAR DO DESENHO PROJETIVORepresentação de vistas como sistema internacional; representação de arestas visíveis e invisíveis; linhas de centro e eixos de simetria.3 - PROJEÇÕES ORTOGONAIS PELO 1º DIEDROPrincípio fundamental; projeções principais; rebatimentos convencionados.4 - PROJEÇÕES ORTOGONAIS PELO 3º DIEDROPrincípio fundamental; projeções principais; rebatimentos convencionados.5 - LEITURA E INTERPRETAÇÃO DE DESENHOSLeitura por meio de esboço em perspectiva e mediante construção de modelos.6 - ESCALASDefinição e normalização7 - DESENHOS COM INSTRUMENTOSRegras para emprego dos esquadros, compasso e régua `"T`"; disposição do desenho nas folhas padronizadas.8 - CORTES E REPRESENTAÇÕES CONVENCIONAISPrincípios fundamentais; aplicações; tipos normalizados; representações e regras para traçado; seções e rupturas.9 - PROJEÇÕES AUXILIARESPrincípios fundamentais; finalidades e aplicações; representações normalizadas.10 - COTAÇÃORegras de colocação e distribuição de cotas.11 - DESENHOS DE CONJUNTOS E DETALHESDefinições; tipos recomendados de legenda e lista de peça; formas de numeração de desenhos; regras práticas para execução e verificação de desenhos.12 - APLICAÇÃO DE TOLERÂNCIAS E AJUSTESDefinição e finalidades; sistema ISO; uso de tabelas e indicação nos desenhos.13 - SÍMBOLOS DE ACABAMENTO SUPERFICIALDefinição; simbologia normalizada; aplicações.14 - DESENHO DE ELEMENTOS DE MÁQUINADefinições, aplicações, tipos, proporções e representações convencionais de: roscas, parafusos, porcas, arruelas, polias, correias e chavetas.15 - DESENHO DE EQUIPAMENTOS E ACESSÓRIOSDesenho de conjunto e detalhes envolvendo elementos de ligação e de máquinas com aplicação de tabelas e catálogos." "1 - INTRODUÇÃO`vApresentação e definição da disciplina, destacando a importância do desenho na engenharia; Normas ABNT e ISO.`v2 - TEORIA ELEMENTAR DO DESENHO PROJETIVO`vRepresentação de vistas como sistema internacional; representação de arestas visíveis e invisíveis; linhas de centro e eixos de simetria.`v3 - PROJEÇÕES ORTOGONAIS PELO 1º DIEDRO`vPrincípio fundamental; projeções principais; rebatimentos convencionados.`v4 - PROJEÇÕES ORTOGONAIS PELO 3º DIEDRO`vPrincípio fundamental; projeções principais; rebatimentos convencionados.`v5 - LEITURA E INTERPRETAÇÃO DE DESENHOS`vLeitura por meio de esboço em perspectiva e mediante construção de modelos.`v6 - ESCALAS`vDefinição e normalização`v7 - DESENHOS COM INSTRUMENTOS`vRegras para emprego dos esquadros, compasso e régua `"T`"; disposição do desenho nas folhas padronizadas.`v8 - CORTES E REPRESENTAÇÕES CONVENCIONAIS`vPrincípios fundamentais; aplicações; tipos normalizados; representações e regras para traçado; seções e rupturas.`v9 - PROJEÇÕES AUXILIARES`vPrincípios fundamentais; finalidades e aplicações; representações normalizadas.`v10 - COTAÇÃO`vRegras de colocação e distribuição de cotas.`v11 - DESENHOS DE CONJUNTOS E DETALHES`vDefinições; tipos recomendados de legenda e lista de peça; formas de numeração de desenhos; regras práticas para execução e verificação de desenhos.`v12 - APLICAÇÃO DE TOLERÂNCIAS E AJUSTES`vDefinição e finalidades; sistema ISO; uso de tabelas e indicação nos desenhos.`v13 - SÍMBOLOS DE ACABAMENTO SUPERFICIAL`vDefinição; simbologia normalizada; aplicações.`v14 - DESENHO DE ELEMENTOS DE MÁQUINA`vDefinições, aplicações, tipos, proporções e representações convencionais de: roscas, parafusos, porcas, arruelas, polias, correias e chavetas.`v15 - DESENHO DE EQUIPAMENTOS E ACESSÓRIOS`vDesenho de conjunto e detalhes envolvendo elementos de ligação e de máquinas com aplicação de tabelas e catálogos."
Write-Output "Edit 3: $ok3"

# Edit 4
$ok4 = Split-WithLineBreaks "- A recuperação deverá consistir de uma prova englobando a matéria toda do semestre.- A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação." "- A recuperação deverá consistir de uma prova englobando a matéria toda do semestre.`v- A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação."
Write-Output "Edit 4: $ok4"

# Edit 5
$ok5 = Split-WithLineBreaks "1 - ABNT - COLETÂNEA DE NORMAS DE DESENHO TÉCNICONormas Técnicas publicadas pela ABNT2 - DESENHO BÁSICO NA ENGENHARIARibeiro, Antonio Clélio - Apostila publicada pela FAENQUIL3 - FUNDAMENTOS DE DIBUJO EM INGENIERIALuzader, Warren J. - Ed. Comp. Editorial Continental - México4 - MANUAL DE DESENHO TÉCNICOManfé, G./ Scarato, G./ Pozza, R. - Ed. Renovada Livros Culturais Ltda.5 - EXPRESSÃO GRÁFICA - DESENHO TÉCNICOHoelsher, R. P./ Springer, C.H./ Dobrovolny, J.S. - Ed. LTC Editora S.A.6 - DESENHO TÉCNICOFrench, Thomas E. - Editora Globo7 - DESENHO TÉCNICOBachmann, A./ Forberg, R - Editora Globo8 - DESENHISTA DE MÁQUINASEscola PRO-TEC" "1 - ABNT - COLETÂNEA DE NORMAS DE DESENHO TÉCNICO`vNormas Técnicas publicadas pela ABNT`v2 - DESENHO BÁSICO NA ENGENHARIA`vRibeiro, Antonio Clélio - Apostila publicada pela FAENQUIL`v3 - FUNDAMENTOS DE DIBUJO EM INGENIERIA`vLuzader, Warren J. - Ed. Comp. Editorial Continental - México`v4 - MANUAL DE DESENHO TÉCNICO`vManfé, G./ Scarato, G./ Pozza, R. - Ed. Renovada Livros Culturais Ltda.`v5 - EXPRESSÃO GRÁFICA - DESENHO TÉCNICO`vHoelsher, R. P./ Springer, C.H./ Dobrovolny, J.S. - Ed. LTC Editora S.A.`v6 - DESENHO TÉCNICO`vFrench, Thomas E. - Editora Globo`v7 - DESENHO TÉCNICO`vBachmann, A./ Forberg, R - Editora Globo`v8 - DESENHISTA DE MÁQUINAS`vEscola PRO-TEC"
Write-Output "Edit 5: $ok5"

